$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix two existing rows: E57 and E58 should become TRUE
$ws.Range("E57").Value = $true
$ws.Range("E58").Value = $true

# Add a new row 61 for "Emma G"
$ws.Range("A61").Value = "Emma G"
$ws.Range("B61").Value = 1
$ws.Range("C61").Value = "Torso"
$ws.Range("D61").Value = $true
$ws.Range("E61").Value = $false
$ws.Range("F61").Value = $false
$ws.Range("G61").Value = $false
$ws.Range("H61").Value = $true

# Update the table range to include the new row
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:I61"))

# Update view state (scroll position + active selection)
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
$ws.Range("E59").Select()
